# Generate Report for Handback
#
# The handback-transform job failed to match the generated handback file
# name (sxmjzysi.fpo) to the expected handoff file name for the
# f38146c0-f1bc-457a-b5d8-dec631394f74 document, for both the zh-cn and
# de-de targets. Reflect that failure on the status report:
#   - flip the "Ready for handoff" status (shared across Overview/zh-cn/de-de)
#     to "Handback transform failed" for that row
#   - record the mismatch detail in the "Error Detail" (column L) cell for
#     that row on both the zh-cn and de-de sheets

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zh = $wb.Worksheets.Item("zh-cn")
$ws_de = $wb.Worksheets.Item("de-de")

# Row 3 on every sheet (Overview!B3/C3, zh-cn!C3, de-de!C3) shares the
# string "Ready for handoff" -- flip every occurrence so the shared text
# changes everywhere it is used.
$ws_overview.Range("B3").Value = "Handback transform failed"
$ws_overview.Range("C3").Value = "Handback transform failed"
$ws_zh.Range("C3").Value = "Handback transform failed"
$ws_de.Range("C3").Value = "Handback transform failed"

# New error-detail cells for the f38146c0... row on the locale sheets.
$ws_zh.Range("L3").Value = "Handback file name: sxmjzysi.fpo is different with handoff file name: f38146c0-f1bc-457a-b5d8-dec631394f74.b01412e4defff509b63db3443706b29e67722fc2.zh-cn."
$ws_de.Range("L3").Value = "Handback file name: sxmjzysi.fpo is different with handoff file name: f38146c0-f1bc-457a-b5d8-dec631394f74.b01412e4defff509b63db3443706b29e67722fc2.de-de."
